$wb = $excel.ActiveWorkbook

# Sheet: ALC (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value2 = 1323.2433
$ws.Range("J17").Value2 = 1323.2433
$ws.Range("L17").Value2 = 3969.7299
$ws.Range("N17").Value2 = -4305.7299
$ws.Range("H28").Value2 = 1610.6471
$ws.Range("I28").Value2 = 1610.6471
$ws.Range("K28").Value2 = 1610.6471
$ws.Range("M28").Value2 = -1125.6471
$ws.Range("H32").Value2 = 2428.1052
$ws.Range("I32").Value2 = 1925.4445
$ws.Range("J32").Value2 = 2880.5
$ws.Range("K32").Value2 = 1925.4445
$ws.Range("L32").Value2 = 2880.5
$ws.Range("M32").Value2 = -1599.4445
$ws.Range("N32").Value2 = -3532.5
$ws.Range("H33").Value2 = 300.35715
$ws.Range("I33").Value2 = 310.42307
$ws.Range("J33").Value2 = 169.5
$ws.Range("K33").Value2 = 310.42307
$ws.Range("L33").Value2 = 169.5
$ws.Range("M33").Value2 = -81.42307
$ws.Range("N33").Value2 = -627.5
$ws.Range("H106").Value2 = 14470.75
$ws.Range("I106").Value2 = 14470.75
$ws.Range("K106").Value2 = 14470.75
$ws.Range("M106").Value2 = -13839.75
$ws.Range("H116").Value2 = 3124
$ws.Range("I116").Value2 = 2999
$ws.Range("J116").Value2 = 3499
$ws.Range("K116").Value2 = 2999
$ws.Range("L116").Value2 = 3499
$ws.Range("M116").Value2 = 443
$ws.Range("N116").Value2 = -10383
$ws.Range("H135").Value2 = 1022.1429
$ws.Range("J135").Value2 = 1069.5
$ws.Range("L135").Value2 = 9625.5
$ws.Range("N135").Value2 = -14695.5
$ws.Range("H137").Value2 = 4367.1
$ws.Range("I137").Value2 = 3945.1667
$ws.Range("K137").Value2 = 11835.5001
$ws.Range("M137").Value2 = -9285.500100000001
$ws.Range("H138").Value2 = 1932.0834
$ws.Range("I138").Value2 = 1932.0834
$ws.Range("K138").Value2 = 5796.2502
$ws.Range("M138").Value2 = -656.2502000000004
$ws.Range("H141").Value2 = 8970.416999999999
$ws.Range("I141").Value2 = 9777.5
$ws.Range("J141").Value2 = 4935
$ws.Range("K141").Value2 = 29332.5
$ws.Range("L141").Value2 = 14805
$ws.Range("M141").Value2 = -24152.5
$ws.Range("N141").Value2 = -25165

# Sheet: ARM (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value2 = 1776.2307
$ws.Range("I32").Value2 = 1691.3948
$ws.Range("K32").Value2 = 1691.3948
$ws.Range("M32").Value2 = -1404.3948
$ws.Range("H45").Value2 = 2449.75
$ws.Range("I45").Value2 = 2449.75
$ws.Range("J45").Value2 = 0
$ws.Range("K45").Value2 = 2449.75
$ws.Range("L45").Value2 = 0
$ws.Range("M45").Value2 = -2072.75
$ws.Range("N45").ClearContents()
$ws.Range("H46").Value2 = 16247.6
$ws.Range("I46").Value2 = 13746
$ws.Range("J46").Value2 = 20000
$ws.Range("K46").Value2 = 13746
$ws.Range("L46").Value2 = 20000
$ws.Range("M46").Value2 = -13427
$ws.Range("N46").Value2 = -20638
$ws.Range("H61").Value2 = 4660.3716
$ws.Range("I61").Value2 = 3681.5925
$ws.Range("K61").Value2 = 3681.5925
$ws.Range("M61").Value2 = -3469.5925
$ws.Range("H74").Value2 = 2663.6667
$ws.Range("I74").Value2 = 1996
$ws.Range("J74").Value2 = 2997.5
$ws.Range("K74").Value2 = 1996
$ws.Range("L74").Value2 = 2997.5
$ws.Range("M74").Value2 = -1122
$ws.Range("N74").Value2 = -4745.5
$ws.Range("H77").Value2 = 2663.6667
$ws.Range("I77").Value2 = 1996
$ws.Range("J77").Value2 = 2997.5
$ws.Range("K77").Value2 = 9980
$ws.Range("L77").Value2 = 14987.5
$ws.Range("M77").Value2 = -5612
$ws.Range("N77").Value2 = -23723.5
$ws.Range("H110").Value2 = 6329.8
$ws.Range("I110").Value2 = 7506
$ws.Range("J110").Value2 = 1625
$ws.Range("K110").Value2 = 7506
$ws.Range("L110").Value2 = 1625
$ws.Range("M110").Value2 = -5461
$ws.Range("N110").Value2 = -5715
$ws.Range("H122").Value2 = 1664
$ws.Range("I122").Value2 = 1664
$ws.Range("K122").Value2 = 4992
$ws.Range("M122").Value2 = -2542
$ws.Range("H132").Value2 = 3139
$ws.Range("I132").Value2 = 3139
$ws.Range("K132").Value2 = 9417
$ws.Range("M132").Value2 = -6887
$ws.Range("H136").Value2 = 4660.3716
$ws.Range("I136").Value2 = 3681.5925
$ws.Range("K136").Value2 = 11044.7775
$ws.Range("M136").Value2 = -8494.7775

# Sheet: BSM (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("H7").Value2 = 29027.5
$ws.Range("I7").Value2 = 0
$ws.Range("K7").Value2 = 0
$ws.Range("M7").ClearContents()
$ws.Range("H94").Value2 = 839.26666
$ws.Range("I94").Value2 = 839.26666
$ws.Range("J94").Value2 = 0
$ws.Range("K94").Value2 = 839.26666
$ws.Range("L94").Value2 = 0
$ws.Range("M94").Value2 = -388.26666
$ws.Range("N94").ClearContents()
$ws.Range("H134").Value2 = 5941.2915
$ws.Range("I134").Value2 = 5561
$ws.Range("K134").Value2 = 16683
$ws.Range("M134").Value2 = -14148

# Sheet: CRP (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value2 = 3334100
$ws.Range("I22").Value2 = 812.625
$ws.Range("J22").Value2 = 10000675
$ws.Range("K22").Value2 = 812.625
$ws.Range("L22").Value2 = 10000675
$ws.Range("M22").Value2 = -462.625
$ws.Range("N22").Value2 = -10001375
$ws.Range("H31").Value2 = 2267
$ws.Range("I31").Value2 = 1863.4
$ws.Range("K31").Value2 = 1863.4
$ws.Range("M31").Value2 = -1568.4
$ws.Range("H34").Value2 = 2267
$ws.Range("I34").Value2 = 1863.4
$ws.Range("K34").Value2 = 1863.4
$ws.Range("M34").Value2 = -1661.4
$ws.Range("H97").Value2 = 0
$ws.Range("J97").Value2 = 0
$ws.Range("L97").Value2 = 0
$ws.Range("N97").ClearContents()
$ws.Range("H105").Value2 = 1872.6666
$ws.Range("I105").Value2 = 1872.6666
$ws.Range("K105").Value2 = 1872.6666
$ws.Range("M105").Value2 = -125.6666
$ws.Range("H107").Value2 = 672.26666
$ws.Range("J107").Value2 = 580.2
$ws.Range("L107").Value2 = 580.2
$ws.Range("N107").Value2 = -4420.2
$ws.Range("H118").Value2 = 94999
$ws.Range("J118").Value2 = 94999
$ws.Range("L118").Value2 = 94999
$ws.Range("N118").Value2 = -98313

# Sheet: CUL (index 5)
$ws = $wb.Worksheets.Item(5)
$ws.Range("H12").Value2 = 44.23077
$ws.Range("J12").Value2 = 48.833332
$ws.Range("L12").Value2 = 146.499996
$ws.Range("N12").Value2 = -492.499996
$ws.Range("H75").Value2 = 449
$ws.Range("I75").Value2 = 130
$ws.Range("J75").Value2 = 528.75
$ws.Range("K75").Value2 = 390
$ws.Range("L75").Value2 = 1586.25
$ws.Range("M75").Value2 = 608
$ws.Range("N75").Value2 = -3582.25
$ws.Range("H78").Value2 = 449
$ws.Range("I78").Value2 = 130
$ws.Range("J78").Value2 = 528.75
$ws.Range("K78").Value2 = 1170
$ws.Range("L78").Value2 = 4758.75
$ws.Range("M78").Value2 = 3822
$ws.Range("N78").Value2 = -14742.75
$ws.Range("H140").Value2 = 558737.6
$ws.Range("I140").Value2 = 669751.8
$ws.Range("K140").Value2 = 2009255.4
$ws.Range("M140").Value2 = -2004075.4

# Sheet: GSM (index 6)
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value2 = 2589.182
$ws.Range("J80").Value2 = 2449.5
$ws.Range("L80").Value2 = 2449.5
$ws.Range("N80").Value2 = -4445.5
$ws.Range("H83").Value2 = 2589.182
$ws.Range("J83").Value2 = 2449.5
$ws.Range("L83").Value2 = 12247.5
$ws.Range("N83").Value2 = -22231.5
$ws.Range("H102").Value2 = 2684.913
$ws.Range("I102").Value2 = 2702.8096
$ws.Range("J102").Value2 = 2497
$ws.Range("K102").Value2 = 2702.8096
$ws.Range("L102").Value2 = 2497
$ws.Range("M102").Value2 = -1080.8096
$ws.Range("N102").Value2 = -5741
$ws.Range("H113").Value2 = 1093.4445
$ws.Range("I113").Value2 = 981.8333
$ws.Range("J113").Value2 = 1316.6666
$ws.Range("K113").Value2 = 981.8333
$ws.Range("L113").Value2 = 1316.6666
$ws.Range("M113").Value2 = 1188.1667
$ws.Range("N113").Value2 = -5656.6666
$ws.Range("H132").Value2 = 4000
$ws.Range("I132").Value2 = 4000
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 12000
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -9470
$ws.Range("N132").ClearContents()

# Sheet: LTW (index 7)
$ws = $wb.Worksheets.Item(7)
$ws.Range("H24").Value2 = 0
$ws.Range("I24").Value2 = 0
$ws.Range("K24").Value2 = 0
$ws.Range("M24").ClearContents()
$ws.Range("H122").Value2 = 3377.7144
$ws.Range("I122").Value2 = 3274
$ws.Range("J122").Value2 = 4000
$ws.Range("K122").Value2 = 9822
$ws.Range("L122").Value2 = 12000
$ws.Range("M122").Value2 = -7372
$ws.Range("N122").Value2 = -16900
$ws.Range("H132").Value2 = 2102.6
$ws.Range("I132").Value2 = 2025.5834
$ws.Range("K132").Value2 = 6076.7502
$ws.Range("M132").Value2 = -3546.7502

# Sheet: WVR (index 8)
$ws = $wb.Worksheets.Item(8)
$ws.Range("H18").Value2 = 1168.2188
$ws.Range("I18").Value2 = 1168.2188
$ws.Range("K18").Value2 = 1168.2188
$ws.Range("M18").Value2 = -995.2188000000001
$ws.Range("H74").Value2 = 19588.8
$ws.Range("I74").Value2 = 18648
$ws.Range("K74").Value2 = 18648
$ws.Range("M74").Value2 = -17712
$ws.Range("H77").Value2 = 19588.8
$ws.Range("I77").Value2 = 18648
$ws.Range("K77").Value2 = 55944
$ws.Range("M77").Value2 = -51264
$ws.Range("H96").Value2 = 3810.375
$ws.Range("I96").Value2 = 3623.5
$ws.Range("J96").Value2 = 3997.25
$ws.Range("K96").Value2 = 3623.5
$ws.Range("L96").Value2 = 3997.25
$ws.Range("M96").Value2 = -2250.5
$ws.Range("N96").Value2 = -6743.25
$ws.Range("H100").Value2 = 1645.8
$ws.Range("I100").Value2 = 427
$ws.Range("K100").Value2 = 854
$ws.Range("M100").Value2 = -313
$ws.Range("H107").Value2 = 708.1539
$ws.Range("I107").Value2 = 693.875
$ws.Range("J107").Value2 = 731
$ws.Range("K107").Value2 = 2081.625
$ws.Range("L107").Value2 = 2193
$ws.Range("M107").Value2 = -161.625
$ws.Range("N107").Value2 = -6033
$ws.Range("H132").Value2 = 6938.3335
$ws.Range("I132").Value2 = 5312
$ws.Range("K132").Value2 = 15936
$ws.Range("M132").Value2 = -13406
